$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# EXERCISE 1 sheet (internal sheet2.xml) — fill in the run count / average
# time for the SPECTRAL NAÏVE (GIGI) row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("EXERCISE 1")

$ws1.Range("C2").Value = 1
$ws1.Range("E2").Value = 82966.5179359912

$ws1.Range("E12").Select()

# ---------------------------------------------------------------------------
# EXERCISE 2 sheet (internal sheet3.xml) — extend the results table with
# per-run columns (LANCIO 1..10) and split the PageRank rows by tolerance.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("EXERCISE 2")

# Make room: one new row for "PAGERANK NAÏVE TOLLERANZA 1E-8" (after HITS
# PARALLEL) and one new row for "PAGERANK VECTORIZED Tolleranza 1e-8" (at
# the very end).
$ws2.Rows("10").Insert()
$ws2.Rows("13").Insert()

# Extend header row with the per-run columns.
$ws2.Range("D1").Value = "LANCIO 1 "
$ws2.Range("E1").Value = "LANCIO 2 "
$ws2.Range("F1").Value = "LANCIO 3 "
$ws2.Range("G1").Value = "LANCIO 4 "
$ws2.Range("H1").Value = "LANCIO 5 "
$ws2.Range("I1").Value = "LANCIO 6 "
$ws2.Range("J1").Value = "LANCIO 7 "
$ws2.Range("K1").Value = "LANCIO 8"
$ws2.Range("L1").Value = "LANCIO 9"
$ws2.Range("M1").Value = "LANCIO 10"

# Row 5 used to be "PAGERANK NAÏVE " -- now split by tolerance and filled
# in with the run data + an average formula.
$ws2.Range("A5").Value = "PAGERANK NAÏVE TOLLERANZA 1E-6"
$ws2.Range("C5").Formula = "=AVERAGE(D5:M5)"
$ws2.Range("D5").Value = 57.1544575691223
$ws2.Range("E5").Value = 57.921657085418701
$ws2.Range("F5").Value = 61.6450967788696
$ws2.Range("G5").Value = 57.5576364994049
$ws2.Range("H5").Value = 56.907742261886597
$ws2.Range("I5").Value = 56.983201742172199
$ws2.Range("J5").Value = 56.932199716567901
$ws2.Range("K5").Value = 57.103562355041497
$ws2.Range("L5").Value = 57.684858083724897
$ws2.Range("M5").Value = 58.887307167053201

# "HITS NAÏVE " row used a stray non-bold style; normalise it like its
# neighbours (all other rows in column A are bold).
$ws2.Range("A6").Font.Bold = $true

# CLOSENESS PARALLEL shifted down to row 11.
$ws2.Range("A11").Value = "CLOSENESS PARALLEL"
$ws2.Range("A11").Font.Bold = $true

# New row: PAGERANK NAÏVE TOLLERANZA 1E-8 (the naive, non-vectorized run).
$ws2.Range("A12").Value = "PAGERANK NAÏVE TOLLERANZA 1E-8"
$ws2.Range("A12").Font.Bold = $true
$ws2.Range("D12").Value = 303.06279754638598

# New row: PAGERANK VECTORIZED Tolleranza 1e-6 (faster, vectorized run).
$ws2.Range("A10").Value = "PAGERANK VECTORIZED Tolleranza 1e-6 "
$ws2.Range("A10").Font.Bold = $true
$ws2.Range("D10").Value = 15.047101974487299

$ws2.Range("A13").Value = "PAGERANK VECTORIZED Tolleranza 1e-8"
$ws2.Range("A13").Font.Bold = $true
$ws2.Range("D13").Value = 13.9894342422485

$ws2.Range("D13").Select()
